$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# Update the employee record value in B2 to match the existing employee ("TestNew")
$ws.Range("B2").Value = "TestNew"

# Move the selection/active cell to B2 (mirrors the saved view state)
$ws.Activate()
$ws.Range("B2").Select()
